{"js": "// Apply the text edits described by the diff to the document body.\n// Each replacement is performed via body.search(...) + insertText(..., replace)\n// so that formatting of the surrounding (unchanged) runs is preserved as much\n// as possible while the edited text is normalised into a single run.\n\nasync function replaceOnce(context, searchText, replacement) {\n  const body = context.document.body;\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"comparing it to 20 other\" -> \"compared to 20 other\"\nawait replaceOnce(\n  context,\n  \"comparing it to 20 other\",\n  \"compared to 20 other\"\n);\n\n// 2. Widest \"range of\" variety removed + \"with\" -> \"having\"/\"projects\"->\"campaigns\"\nawait replaceOnce(\n  context,\n  \"from the widest range of variety for categories \u2013 9 in total. Even with having most number of projects and variety of category, only over a half of the projects were\",\n  \"from the widest variety of categories \u2013 9 in total. Even having most number of projects and variety of category, only over a half of the campaigns were\"\n);\n\n// 3. Insert new sentence about 11 out of 21 countries.\nawait replaceOnce(\n  context,\n  \"the sample data given. While this is true\",\n  \"the sample data given. 11 out of the 21 countries have Theater as the top category for with the greatest number of projects launched. While this is true\"\n);\n\n// 4. Large rewrite of the \"limitations\" / trends paragraph.\nawait replaceOnce(\n  context,\n  \"The sample data does not contain projects from the last two years. Trends could easily be affected by technology which has changed a lot quickly over just the past couple of years. There has been many changes and improvements, even with laws, when it comes to data visibility/availability as well as security that has created a huge effect on how products and services are marketed to the public.\",\n  \"The sample data does not contain projects launched in the last two years. Kickstarter could have a better position in the market during this time, therefore making the data more useful for analysis. Trends is also easily affected by technology which has changed a lot quickly in the last couple of years. There has been many changes and improvements, even with laws, when it comes to data visibility/availability as well as security that has affected markets of all kinds in different ways.\"\n);\n\n// 5. \"for analysis to determine\" -> \"for data analysts to determine\"\nawait replaceOnce(\n  context,\n  \"this could be a useful information for analysis to determine whether\",\n  \"this could be a useful information for data analysts to determine whether\"\n);\n\n// 6. \"each project.\" -> \"each of the campaigns.\"\nawait replaceOnce(\n  context,\n  \"get audience for each project. It is important to reach\",\n  \"get audience for each of the campaigns. It is important to reach\"\n);\n\n// 7. \"Familiarity of different cultures and countries to Kickstarter technologies\"\nawait replaceOnce(\n  context,\n  \"Familiarity of different cultures and countries to Kickstarter technologies\",\n  \"Access and familiarity of different cultures to the service and technologies used by Kickstarter\"\n);\n\n// 8. \"who are not familiar with this...\" rewrite + new trailing sentence.\nawait replaceOnce(\n  context,\n  \"who are not familiar with this can see this as a factor to not have confidence on the service. It could be one reason why Kickstarter is used most in the US than in any country in the dataset.\",\n  \"who are not familiar with the technology can see this as a reason to not use the service due to lack of confidence. It could be one reason why Kickstarter is used most in the US than in any country in the dataset. Organizations from other parts of the world could also not have access to these services contributing to less projects launched from those countries.\"\n);\n\n// 9. \"bigger amounts could be something intimidating\" -> \"... be seen as something intimidating\"\nawait replaceOnce(\n  context,\n  \"bigger amounts could be something intimidating\",\n  \"bigger amounts could be seen as something intimidating\"\n);\n", "ps1": "# Applies the text edits described by the diff using Word COM Find/Replace.\n# wdReplaceAll = 2, wdFindContinue = 1\n\n$d = $word.ActiveDocument\n$enDash = [char]0x2013\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# 1. \"comparing it to 20 other\" -> \"compared to 20 other\"\nReplace-Text \"comparing it to 20 other\" \"compared to 20 other\"\n\n# 2. Widest \"range of\" variety removed + \"with\" -> \"having\"/\"projects\"->\"campaigns\"\n$find2 = \"from the widest range of variety for categories \" + $enDash + \" 9 in total. Even with having most number of projects and variety of category, only over a half of the projects were\"\n$repl2 = \"from the widest variety of categories \" + $enDash + \" 9 in total. Even having most number of projects and variety of category, only over a half of the campaigns were\"\nReplace-Text $find2 $repl2\n\n# 3. Insert new sentence about 11 out of 21 countries.\nReplace-Text \"the sample data given. While this is true\" \"the sample data given. 11 out of the 21 countries have Theater as the top category for with the greatest number of projects launched. While this is true\"\n\n# 4. Large rewrite of the \"limitations\" / trends paragraph.\n$find4 = \"The sample data does not contain projects from the last two years. Trends could easily be affected by technology which has changed a lot quickly over just the past couple of years. There has been many changes and improvements, even with laws, when it comes to data visibility/availability as well as security that has created a huge effect on how products and services are marketed to the public.\"\n$repl4 = \"The sample data does not contain projects launched in the last two years. Kickstarter could have a better position in the market during this time, therefore making the data more useful for analysis. Trends is also easily affected by technology which has changed a lot quickly in the last couple of years. There has been many changes and improvements, even with laws, when it comes to data visibility/availability as well as security that has affected markets of all kinds in different ways.\"\nReplace-Text $find4 $repl4\n\n# 5. \"for analysis to determine\" -> \"for data analysts to determine\"\nReplace-Text \"this could be a useful information for analysis to determine whether\" \"this could be a useful information for data analysts to determine whether\"\n\n# 6. \"each project.\" -> \"each of the campaigns.\"\nReplace-Text \"get audience for each project. It is important to reach\" \"get audience for each of the campaigns. It is important to reach\"\n\n# 7. \"Familiarity of different cultures and countries to Kickstarter technologies\"\nReplace-Text \"Familiarity of different cultures and countries to Kickstarter technologies\" \"Access and familiarity of different cultures to the service and technologies used by Kickstarter\"\n\n# 8. \"who are not familiar with this...\" rewrite + new trailing sentence.\n$find8 = \"who are not familiar with this can see this as a factor to not have confidence on the service. It could be one reason why Kickstarter is used most in the US than in any country in the dataset.\"\n$repl8 = \"who are not familiar with the technology can see this as a reason to not use the service due to lack of confidence. It could be one reason why Kickstarter is used most in the US than in any country in the dataset. Organizations from other parts of the world could also not have access to these services contributing to less projects launched from those countries.\"\nReplace-Text $find8 $repl8\n\n# 9. \"bigger amounts could be something intimidating\" -> \"... be seen as something intimidating\"\nReplace-Text \"bigger amounts could be something intimidating\" \"bigger amounts could be seen as something intimidating\"\n"}
